$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet1 ("Sheet1"): split the old Trt codes (YF/YM/OF) into more
# granular codes (YFD/YMT/YMD/OFD) on a per-cage basis ---
$ws1.Cells.Item(2, 3).Value = "YFD"
$ws1.Cells.Item(3, 3).Value = "YFD"
$ws1.Cells.Item(4, 3).Value = "YFD"
$ws1.Cells.Item(5, 3).Value = "YFD"
$ws1.Cells.Item(6, 3).Value = "YMT"
$ws1.Cells.Item(7, 3).Value = "YMT"
$ws1.Cells.Item(8, 3).Value = "YMT"
$ws1.Cells.Item(9, 3).Value = "YMD"
$ws1.Cells.Item(10, 3).Value = "YMD"
$ws1.Cells.Item(11, 3).Value = "YMD"
$ws1.Cells.Item(12, 3).Value = "YMD"
$ws1.Cells.Item(13, 3).Value = "YMD"
$ws1.Cells.Item(14, 3).Value = "YMD"
$ws1.Cells.Item(15, 3).Value = "YMT"
$ws1.Cells.Item(16, 3).Value = "YMT"
$ws1.Cells.Item(17, 3).Value = "YMT"
$ws1.Cells.Item(18, 3).Value = "YFD"
$ws1.Cells.Item(19, 3).Value = "YFD"
$ws1.Cells.Item(20, 3).Value = "YMT"
$ws1.Cells.Item(21, 3).Value = "YMT"
$ws1.Cells.Item(22, 3).Value = "YMT"
$ws1.Cells.Item(23, 3).Value = "OFD"
$ws1.Cells.Item(24, 3).Value = "OFD"
$ws1.Cells.Item(25, 3).Value = "OFD"
$ws1.Cells.Item(26, 3).Value = "OFD"
$ws1.Cells.Item(27, 3).Value = "OFD"
$ws1.Cells.Item(28, 3).Value = "OFD"

# --- Sheet3 ("Sheet3"): add a new "Trt" column (E) with the same
# granular treatment codes, matching each cage's value on Sheet1 ---
$ws3.Cells.Item(1, 5).Value = "Trt"
$ws3.Cells.Item(2, 5).Value = "YFD"
$ws3.Cells.Item(3, 5).Value = "YFD"
$ws3.Cells.Item(4, 5).Value = "YMT"
$ws3.Cells.Item(5, 5).Value = "YMD"
$ws3.Cells.Item(6, 5).Value = "OFD"
$ws3.Cells.Item(7, 5).Value = "OFD"
$ws3.Cells.Item(8, 5).Value = "YMD"
$ws3.Cells.Item(9, 5).Value = "YMD"
$ws3.Cells.Item(10, 5).Value = "YMT"
$ws3.Cells.Item(11, 5).Value = "YFD"
$ws3.Cells.Item(12, 5).Value = "YMT"
$ws3.Cells.Item(13, 5).Value = "OFD"

# Match header style (s="6") of C1/D1 for the new E1 header cell
$ws3.Range("D1").Copy()
$ws3.Range("E1").PasteSpecial(-4122)

# --- Selection / active-sheet bookkeeping: previously Sheet1 had the
# selected cell/tab, now Sheet3 is the active tab with its own
# selection, and Sheet1 keeps a plain (non-active) selection ---
$ws1.Range("C29").Select()
$ws3.Activate()
$ws3.Range("E14").Select()
